$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row below the first data row (old row 2) for the
# "Java interview progarms" (typo) variant question. This shifts the
# former rows 3-10 down to rows 4-11.
$ws.Rows(3).Insert()

# Copy formatting from the row above into the freshly inserted (blank) row
# so the new row matches the look (style, height) of its neighbours.
$ws.Range("A2:D2").Copy()
$ws.Range("A3:D3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the newly inserted row 3 (serial no. 2 in the sheet's own numbering)
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Y"
$ws.Range("C3").Value = "gpt-4.1-mini"
$ws.Range("D3").Value = "Create Java interview progarms covering Core Java, OOPs, Collections, Exceptions, Java 8, and Multithreading.  Note: Always provide the question in serial number format"

# Renumber the S.No column (A) for the rows that followed, since they
# shifted down by one row but keep their original sequential numbers.
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7
$ws.Range("A9").Value = 8
$ws.Range("A10").Value = 9

# Update the question text (column D) for the shifted rows to append the
# new trailing note used across all rows.
$ws.Range("D4").Value = "Create Selenium WebDriver interview questions covering locators, waits, alerts, frames, windows, Actions class, JavaScriptExecutor, POM, and Selenium Grid.  Note: Always provide the question in serial number format"
$ws.Range("D5").Value = "Create TestNG interview questions covering annotations, assertions, DataProvider, testng.xml, parallel execution, and listeners. Note: Always provide the question in serial number format"
$ws.Range("D6").Value = "Create Maven interview questions covering pom.xml structure, dependencies, build lifecycle, profiles, and Surefire/Failsafe plugins.  Note: Always provide the question in serial number format"
$ws.Range("D7").Value = "Create API testing interview questions covering REST principles, HTTP methods, status codes, authentication, request/response validation, and Rest Assured.  Note: Always provide the question in serial number format"
$ws.Range("D8").Value = "Include basic, intermediate, and advanced level questions for each technology.  Note: Always provide the question in serial number format"
$ws.Range("D9").Value = "Mix theoretical, scenario-based, and practical questions.Focus on real-world automation framework design and best practices.  Note: Always provide the question in serial number format"

# Row 10 keeps the "Target candidates..." text (previously row 9) and row 11
# (brand-new row appended at the bottom) gets the "Do not include
# explanations..." text that used to live in row 10.
$ws.Range("D10").Value = "Target candidates with 3–8 years of QA automation experience.Output questions in a clean, line-by-line, serial-numbered format.  Note: Always provide the question in serial number format"

# Also fix up the also-shifted original row 2 question text (append the note).
$ws.Range("D2").Value = "Create Java interview questions covering Core Java, OOPs, Collections, Exceptions, Java 8, and Multithreading. Note: Always provide the question in serial number format"

# Append brand-new row 11 at the bottom of the table.
$ws.Range("A10:D10").Copy()
$ws.Range("A11:D11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Y"
$ws.Range("C11").Value = "gpt-4.1-mini"
$ws.Range("D11").Value = "Do not include explanations—only questions.This format is concise, structured, and LLM-friendly.  Note: Always provide the question in serial number format"

# Explicitly set row heights to match the new (longer) wrapped text in
# column D, mirroring how Excel would auto-fit these rows after the text
# change.
$ws.Rows(2).RowHeight = 41.4
$ws.Rows(3).RowHeight = 41.4
$ws.Rows(4).RowHeight = 55.2
$ws.Rows(5).RowHeight = 41.4
$ws.Rows(6).RowHeight = 41.4
$ws.Rows(7).RowHeight = 55.2
$ws.Rows(8).RowHeight = 41.4
$ws.Rows(9).RowHeight = 41.4
$ws.Rows(10).RowHeight = 55.2
$ws.Rows(11).RowHeight = 41.4
